$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "RideDate"
$ws.Range("B1").Value = "MileageStart"
$ws.Range("C1").Value = "MileageEnd"
$ws.Range("D1").Value = "Fill Up"
$ws.Range("E1").Value = "Gallons"
$ws.Range("F1").Value = "PricePerGallon"
$ws.Range("G1").Value = "RideRoute"
$ws.Range("H1").Value = "RideDescription"
$ws.Range("I1").Value = "ImagePath"

# Row 2 - initial fill up / starting ride
$ws.Range("A2").Value = 43972
$ws.Range("B2").Value = 415
$ws.Range("C2").Value = 415
$ws.Range("D2").Value = "Yes"
$ws.Range("E2").Value = 2.2
$ws.Range("F2").Value = 1.95
$ws.Range("G2").Value = "Gas Station"
$ws.Range("H2").Value = "Initial Fill Up and Starting Ride"

# Row 3 - ride to Nielsens
$ws.Range("A3").Value = 43972
$ws.Range("B3").Value = 415
$ws.Range("C3").Value = 485
$ws.Range("D3").Value = "No"
$ws.Range("G3").Value = "Down to Nielsens, Lake Villa, IL"

# Date formatting for RideDate column entries
$ws.Range("A2:A3").NumberFormat = "m/d/yyyy"

# Column widths
$ws.Columns.Item(1).ColumnWidth = 13.28515625
$ws.Columns.Item(2).ColumnWidth = 15
$ws.Columns.Item(3).ColumnWidth = 15.140625
$ws.Columns.Item(6).ColumnWidth = 18.42578125
$ws.Columns.Item(7).ColumnWidth = 29.28515625
$ws.Columns.Item(8).ColumnWidth = 28.85546875
$ws.Columns.Item(9).ColumnWidth = 14

# Selection matches final saved state
$ws.Range("G4").Select()
